$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.1619
$ws.Range("A10").Value = -20.44259999999997
$ws.Range("A12").Value = -22.51650000000004
$ws.Range("B13").Value = 5.873499999999996
$ws.Range("A18").Value = -22.49290000000003
$ws.Range("C20").Value = -14.92839999999999
